$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing value for week 15 (row 16)
$ws.Range("B16").Value = 425

# Update existing value for week 36 (row 37)
$ws.Range("B37").Value = 449

# Add new row for week 37 (row 38)
$ws.Range("A38").Value = 37
$ws.Range("B38").Value = 503

# Add new row for week 38 (row 39)
$ws.Range("A39").Value = 38
$ws.Range("B39").Value = 1
